{"js": "// Load all body paragraphs so we can inspect their text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// 1) Delete the standalone intro paragraph entirely.\n// 2) Prefix the three material/item lines with \"\u2022 \".\nconst toDelete = [];\nconst toPrefix = [];\n\nfor (const p of paragraphs.items) {\n  const t = p.text.trim();\n  if (t === \"Cateteres diagn\u00f3sticos decapolar e quadripolares.\") {\n    toDelete.push(p);\n  } else if (\n    t === \"Cateter Decapolar \u2013 EasyFinder\u2122 Deca + conector\" ||\n    t === \"Cateter Quadripolar \u2013 EasyFinder\u2122 Quad + conector (2x)\" ||\n    t === \"Introdutor \u2013 3\"\n  ) {\n    toPrefix.push(p);\n  }\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nfor (const p of toPrefix) {\n  p.insertText(\"\u2022 \", \"Start\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$deleteText = \"Cateteres diagn\u00f3sticos decapolar e quadripolares.\"\n$prefixTexts = @(\n    \"Cateter Decapolar \u2013 EasyFinder\u2122 Deca + conector\",\n    \"Cateter Quadripolar \u2013 EasyFinder\u2122 Quad + conector (2x)\",\n    \"Introdutor \u2013 3\"\n)\n\n# Walk paragraphs back-to-front so deleting one doesn't disturb the\n# indices of paragraphs we still need to visit.\n$count = $d.Paragraphs.Count\nfor ($i = $count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text.TrimEnd([char]13, [char]7)\n\n    if ($t -eq $deleteText) {\n        $p.Range.Delete()\n    }\n    elseif ($prefixTexts -contains $t) {\n        $p.Range.InsertBefore(\"\u2022 \")\n    }\n}\n"}
